$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '50.134.40'
$ws.Range("E2").Value = '  +4.20%  '
$ws.Range("D3").Value = '2.659.94'
$ws.Range("E3").Value = '  +6.62%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''114.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +8.30%  '
$ws.Range("D6").Value = '''326.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.79%  '
$ws.Range("E7").Value = '  +2.08%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '''0.558'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.54%  '
$ws.Range("D10").Value = '''41.38'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.26%  '
$ws.Range("D11").Value = '''20.15'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.58%  '
$ws.Range("E12").Value = '  +3.00%  '
$ws.Range("E13").Value = '  +0.38%  '
$ws.Range("D14").Value = '''7.41'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.33%  '
$ws.Range("D15").Value = '3.073.73'
$ws.Range("E15").Value = '  +6.39%  '
$ws.Range("D16").Value = '2.658.50'
$ws.Range("E16").Value = '  +6.32%  '
$ws.Range("D17").Value = '''0.878'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.87%  '
$ws.Range("D18").Value = '50.055.02'
$ws.Range("E18").Value = '  +4.19%  '
$ws.Range("E19").Value = '  +3.62%  '
$ws.Range("D20").Value = '''6.79'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.16%  '
$ws.Range("D21").Value = '''2.93'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.49%  '
$ws.Range("D22").Value = '0.0₃0962'
$ws.Range("E22").Value = '  +3.45%  '
$ws.Range("D23").Value = '''72.57'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.11%  '
$ws.Range("D24").Value = '''276.93'
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = '  +3.13%  '
$ws.Range("D26").Value = '''27.02'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.94%  '
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").Value = '''10.07'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.37%  '
$ws.Range("D29").Value = '''36.91'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.62%  '
$ws.Range("E30").Value = '  -2.19%  '
$ws.Range("D31").Value = '''0.143'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.96%  '
$ws.Range("D32").Value = '''50.21'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.84%  '
$ws.Range("E33").Value = '  +4.07%  '
$ws.Range("D34").Value = '''19.77'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.58%  '
$ws.Range("E35").Value = '  +5.82%  '
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("E37").Value = '  +9.18%  '
$ws.Range("E38").Value = '  +6.77%  '
$ws.Range("E39").Value = '  +8.75%  '
$ws.Range("D40").Value = '''0.114'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.36%  '
$ws.Range("D41").Value = '''124.04'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.97%  '
$ws.Range("D42").Value = '''22.52'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.35%  '
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("E44").Value = '  +5.46%  '
$ws.Range("D45").Value = '2.108.66'
$ws.Range("E45").Value = '  +5.38%  '
$ws.Range("D46").Value = '''3.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.79%  '
$ws.Range("D47").Value = '''2.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +13.30%  '
$ws.Range("D48").Value = '''2.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.44%  '
$ws.Range("E49").Value = '  +2.20%  '
$ws.Range("D50").Value = '''5.38'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.44%  '
$ws.Range("D51").Value = '''60.24'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.44%  '
